$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.257.75"
$ws.Range("E2").Value = "  +4.65%  "
$ws.Range("D3").Value = "3.246.36"
$ws.Range("E3").Value = "  +2.33%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'577.60"
$ws.Range("E5").Value = "  +2.45%  "
$ws.Range("D6").Value = "'178.73"
$ws.Range("E6").Value = "  +6.11%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'0.600"
$ws.Range("E8").Value = "  -3.31%  "
$ws.Range("D9").Value = "3.247.99"
$ws.Range("E10").Value = "  +4.56%  "
$ws.Range("D11").Value = "'6.75"
$ws.Range("E11").Value = "  +2.63%  "
$ws.Range("D12").Value = "'0.413"
$ws.Range("E12").Value = "  +4.44%  "
$ws.Range("D13").Value = "3.811.73"
$ws.Range("E13").Value = "  +2.49%  "
$ws.Range("E14").Value = "  +0.69%  "
$ws.Range("D15").Value = "'27.92"
$ws.Range("E15").Value = "  +2.20%  "
$ws.Range("D16").Value = "67.203.33"
$ws.Range("E16").Value = "  +4.60%  "
$ws.Range("D17").Value = "'0.0000168"
$ws.Range("E17").Value = "  +3.01%  "
$ws.Range("D18").Value = "3.249.90"
$ws.Range("E18").Value = "  +2.33%  "
$ws.Range("D19").Value = "'5.82"
$ws.Range("E19").Value = "  +1.72%  "
$ws.Range("D20").Value = "'13.36"
$ws.Range("E20").Value = "  +3.25%  "
$ws.Range("D21").Value = "'373.80"
$ws.Range("E21").Value = "  +5.92%  "
$ws.Range("D22").Value = "'7.57"
$ws.Range("E22").Value = "  +5.28%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "'71.29"
$ws.Range("E24").Value = "  +4.01%  "
$ws.Range("E25").Value = "  +1.14%  "
$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").Value = "3.387.40"
$ws.Range("E26").Value = "  +2.13%  "
$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D27").Value = "'0.0000119"
$ws.Range("E27").Value = "  +2.49%  "
$ws.Range("D28").Value = "'9.70"
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("E29").Value = "  +2.47%  "
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("E31").Value = "  +4.25%  "
$ws.Range("D32").Value = "'5.62"
$ws.Range("E32").Value = "  +1.43%  "
$ws.Range("D33").Value = "'22.58"
$ws.Range("E33").Value = "  +2.64%  "
$ws.Range("E34").Value = "  +7.31%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").Value = "'6.81"
$ws.Range("E36").Value = "  +2.85%  "
$ws.Range("D37").Value = "'162.88"
$ws.Range("E37").Value = "  +6.37%  "
$ws.Range("D38").Value = "'1.49"
$ws.Range("E38").Value = "  +4.41%  "
$ws.Range("D39").Value = "'0.862"
$ws.Range("E39").Value = "  +5.61%  "
$ws.Range("D40").Value = "'1.85"
$ws.Range("E40").Value = "  +9.61%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'6.81"
$ws.Range("E41").Value = "  +13.28%  "
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").Value = "'26.74"
$ws.Range("E42").Value = "  +1.87%  "
$ws.Range("D43").Value = "'2.61"
$ws.Range("E43").Value = "  +6.39%  "
$ws.Range("D44").Value = "'364.73"
$ws.Range("E44").Value = "  +14.89%  "
$ws.Range("D45").Value = "2.766.25"
$ws.Range("E45").Value = "  +5.78%  "
$ws.Range("E46").Value = "  +5.32%  "
$ws.Range("D47").Value = "'26.01"
$ws.Range("E47").Value = "  +9.76%  "
$ws.Range("D48").Value = "'40.41"
$ws.Range("E48").Value = "  +2.84%  "
$ws.Range("E49").Value = "  +5.03%  "
$ws.Range("E50").Value = "  +3.16%  "
$ws.Range("E51").Value = "  +0.84%  "
